# CST-452 Sprint Burn Down - advance from Week Three to Week Four.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: Sprint number ---
$ws.Range("B2").Value = 6

# --- Week dates (row 3) shift forward by 14 days (one Monday + one week) ---
$ws.Range("G3").Value = 44641
$ws.Range("H3").Value = 44642
$ws.Range("I3").Value = 44643
$ws.Range("J3").Value = 44644
$ws.Range("K3").Value = 44645
$ws.Range("L3").Value = 44646
$ws.Range("M3").Value = 44647

# --- Row 5: Rework Form Styles / N/A / Adjusted Login for to better fit the application styles. / Marc ---
$ws.Range("B5").Value = "N/A"
$ws.Range("C5").Value = "Rework Form Styles"
$ws.Range("D5").Value = "Adjusted Login for to better fit the application styles."
$ws.Range("E5").Value = "Marc"
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0

# --- Row 6: (blank user story / task id) / Changed Registration form to match Login form. / " ---
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "Changed Registration form to match Login form."
$ws.Range("E6").Value = """"
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Rows(6).RowHeight = 43.7

# --- Row 7: N/A / Reviewed Application Documentation. / Updated documentations. / " ---
$ws.Range("B7").Value = "N/A"
$ws.Range("C7").Value = "Reviewed Application Documentation."
$ws.Range("D7").Value = "Updated documentations."
$ws.Range("E7").Value = """"
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 2

# --- Row 8: N/A / Walkthrough Presentation. / Recorded Final Presentation. / " ---
# Leading apostrophe keeps the cell's quote-prefix (forced-text) style intact.
$ws.Range("B8").Value = "'N/A"
$ws.Range("C8").Value = "Walkthrough Presentation."
$ws.Range("D8").Value = "Recorded Final Presentation."
$ws.Range("E8").Value = """"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 1

# --- Row 9: cleared out (previously Application Deployment task) ---
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("F9").Value = 0
$ws.Range("K9").Value = 0
$ws.Rows(9).RowHeight = 18

# --- Row 10: cleared out (previously deployment DAO task) ---
$ws.Range("D10").Value = ""
$ws.Range("F10").Value = 0
$ws.Range("L10").Value = 0
$ws.Rows(10).RowHeight = 18

# --- View: selection moved, zoomed scale changed ---
$ws.Range("L10").Select()
$ws.PageSetup.Zoom = 71
